# Update app with plotly dependency
# Fills in newly-available Destatis/HWWI data for 2025 Q2 (row 43) and
# 2025 Q3 (row 44), and removes the still-empty 2025 Q4 placeholder row
# (old row 45) now that it is no longer the "next" quarter to fill in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 2025 Q2 (row 43): price-index columns C:F and the volume columns
# O:W were still blank before this update -> fill them with the newly
# published figures.
$ws.Range("C43").Value = 206.45
$ws.Range("D43").Value = 152.44
$ws.Range("E43").Value = 155.37
$ws.Range("F43").Value = 376.9

$ws.Range("O43").Value = 509
$ws.Range("P43").Value = 90306.666666666686
$ws.Range("Q43").Value = 34012667
$ws.Range("R43").Value = 1178328666
$ws.Range("S43").Value = 6355352000
$ws.Range("T43").Value = 3367534333
$ws.Range("U43").Value = 2987817666
$ws.Range("V43").Value = 1757211000
$ws.Range("W43").Value = 1230606667

# ---- 2025 Q3 (row 44): the Konjunktur/HWWI survey indicators I:N were
# still blank -> fill them in as well.
$ws.Range("I44").Value = -3.33
$ws.Range("J44").Value = -4.45
$ws.Range("K44").Value = -26.670000000000009
$ws.Range("L44").Value = -18.39
$ws.Range("M44").Value = 3.33
$ws.Range("N44").Value = -40

# ---- 2025 Q4 (old row 45) had no data at all yet -> drop the row so the
# sheet again ends right after the last populated quarter.
$ws.Rows(45).Delete()

# ---- Match the author's on-save selection/scroll state: the active
# cell sits on the freshly-entered block of row 44.
$ws.Range("O44:W44").Select()

# ---- Page margins were reset to Excel's normal defaults (Normal
# template: 0.7/0.7/0.75/0.75/0.3/0.3 in) instead of the oversized
# custom margins the sheet previously had.
$ws.PageSetup.LeftMargin = 50.4
$ws.PageSetup.RightMargin = 50.4
$ws.PageSetup.TopMargin = 54
$ws.PageSetup.BottomMargin = 54
$ws.PageSetup.HeaderMargin = 21.599999999999998
$ws.PageSetup.FooterMargin = 21.599999999999998
